# "data prep test run" - reorder the node rows on the "Definition" and
# "Nodes" sheets into a new (shuffled) order. The row for each named
# object keeps all of its original column data; only the row position
# changes.

$wb = $excel.ActiveWorkbook

# Desired new order (by Object_Name) for both the "Definition" sheet
# (rows 8-20) and the "Nodes" sheet (rows 2-14).
$newOrder = @(
    "Vaporized_Carbon_Dioxide",
    "E-Methanol_storage_Kasso",
    "Power_Wholesale",
    "Waste_Heat",
    "E-Methanol_Kasso",
    "Hydrogen_storage_Kasso",
    "Power_Kasso",
    "Water",
    "District_Heating",
    "Hydrogen_Kasso",
    "Carbon_Dioxide",
    "Raw_Methanol",
    "Steam"
)

### --- "Definition" sheet: rows 8-20, column A (Object_Name) only --- ###

$wsDef = $wb.Worksheets.Item("Definition")
$defFirstRow = 8

# Snapshot the existing rows 8-20 (Object_Name + Category) keyed by name.
$defByName = @{}
for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $r = $defFirstRow + $i
    $name = $wsDef.Cells.Item($r, 1).Value2
    $category = $wsDef.Cells.Item($r, 2).Value2
    $defByName[$name] = $category
}

# Write the rows back out in the new order. Text values are prefixed with a
# leading apostrophe so the engine stores them as plain text instead of
# auto-coercing look-alike values (e.g. "true"/"false") into booleans; the
# cell style is then reset so the apostrophe doesn't leave a lingering
# "quote prefix" style behind.
for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $r = $defFirstRow + $i
    $name = $newOrder[$i]

    $cellA = $wsDef.Cells.Item($r, 1)
    $cellA.Value = "'" + $name
    $cellA.Style = "Normal"

    $cellB = $wsDef.Cells.Item($r, 2)
    $cellB.Value = "'" + $defByName[$name]
    $cellB.Style = "Normal"
}

### --- "Nodes" sheet: rows 2-14, columns A-G --- ###

$wsNodes = $wb.Worksheets.Item("Nodes")
$nodesFirstRow = 2
$lastCol = 7

# Columns B (Category), C (balance_type) and D (has_state) hold text
# values; E/F/G (node_state_cap / frac_state_loss / node_slack_penalty)
# hold numbers.
$textCols = @(2, 3, 4)

# Snapshot the existing rows 2-14 (all columns) keyed by Object_Name.
$nodesByName = @{}
for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $r = $nodesFirstRow + $i
    $name = $wsNodes.Cells.Item($r, 1).Value2
    $rowData = @()
    for ($c = 2; $c -le $lastCol; $c++) {
        $rowData += , ($wsNodes.Cells.Item($r, $c).Value2)
    }
    $nodesByName[$name] = $rowData
}

# Write the rows back out in the new order. Text values are prefixed with a
# leading apostrophe so the engine stores them as plain text instead of
# auto-coercing look-alike values (e.g. "true"/"false") into booleans.
for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $r = $nodesFirstRow + $i
    $name = $newOrder[$i]
    $rowData = $nodesByName[$name]

    $cellA = $wsNodes.Cells.Item($r, 1)
    $cellA.Value = "'" + $name
    $cellA.Style = "Normal"

    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $wsNodes.Cells.Item($r, $c)
        $val = $rowData[$c - 2]
        if (($null -eq $val) -or ($val -eq "")) {
            # Leave the cell blank, matching the original empty cell.
            $cell.Value = ""
        } elseif ($textCols -contains $c) {
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
